# Tutorial 6 solution update: dates in column A change from dd/mm/yyyy to
# dd-mm-yyyy text, and the first attendance record (rows 3-4) flips from
# "Invalid" to "Real" (D and G columns go from 0 to 1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateUpdates = @(
    @{ Row = 3;  Date = "28-07-2022" },
    @{ Row = 4;  Date = "01-08-2022" },
    @{ Row = 5;  Date = "04-08-2022" },
    @{ Row = 6;  Date = "08-08-2022" },
    @{ Row = 7;  Date = "11-08-2022" },
    @{ Row = 8;  Date = "15-08-2022" },
    @{ Row = 9;  Date = "18-08-2022" },
    @{ Row = 10; Date = "22-08-2022" },
    @{ Row = 11; Date = "25-08-2022" },
    @{ Row = 12; Date = "29-08-2022" },
    @{ Row = 13; Date = "01-09-2022" },
    @{ Row = 14; Date = "05-09-2022" },
    @{ Row = 15; Date = "08-09-2022" },
    @{ Row = 16; Date = "12-09-2022" },
    @{ Row = 17; Date = "15-09-2022" },
    @{ Row = 18; Date = "19-09-2022" },
    @{ Row = 19; Date = "22-09-2022" },
    @{ Row = 20; Date = "26-09-2022" },
    @{ Row = 21; Date = "29-09-2022" }
)

foreach ($u in $dateUpdates) {
    $cell = $ws.Cells.Item($u.Row, 1)
    # Force text so Excel's auto date-recognition doesn't turn the
    # dd-mm-yyyy string into a date serial, then drop the temporary
    # number-format override so the cell keeps its original (default) style.
    $cell.NumberFormat = "@"
    $cell.Value = $u.Date
    $cell.ClearFormats()
}

# Row 3 and row 4 attendance records move from "Invalid" to "Real".
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("G4").Value = 1
